# A02RSDS Pricing Report - refreshed data pull.
# The report was re-generated a week later (06/13/2024 -> 06/20/2024), which
# shifts the "Generated:"/"Pricing Date:" footer stamped on every section of
# the report and updates the underlying cost/labor-hour figures that came
# out of the new pricing run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Generated: <date>, <time>" footer, repeated once per report section ---
$generatedStamp = "Generated: 06/20/2024, 04:01 PM"
foreach ($cellRef in @("A3","A58","A101","A136","A171","A206","A238","A261","A285","A294","A358")) {
    $ws.Range($cellRef).Value = $generatedStamp
}

# --- "Pricing Date: <date>" footer, repeated once per report section ---
$pricingDateStamp = "Pricing Date: 06/20/2024"
foreach ($cellRef in @("A10","A65","A108","A143","A178","A213","A244","A267","A291","A301","A364")) {
    $ws.Range($cellRef).Value = $pricingDateStamp
}

# --- Part A / FRC: 845C construction-details block (rows 12-56) ---
$ws.Range("D28").Value = 264.048
$ws.Range("E28").Value = 5541
$ws.Range("G28").Value = 2770.5
$ws.Range("D33").Value = 1945.208
$ws.Range("G33").Value = 30035.01
$ws.Range("E36").Value = 39250.36
$ws.Range("G36").Value = 4121.29
$ws.Range("G38").Value = 43371.65
$ws.Range("G43").Value = 238350.26
$ws.Range("G46").Value = 238430.26
$ws.Range("G47").Value = 30045.68
$ws.Range("G49").Value = 30045.68
$ws.Range("G50").Value = 353432.94
$ws.Range("G53").Value = 439606.94
$ws.Range("E56").Value = 273.39          # LABOR HOURS: (grand total)

# --- FRC: 822C block (rows 90-99) ---
$ws.Range("G90").Value = 1308.1
$ws.Range("E90").Value = "11.90 HRS"     # SPLICING LABOR hours label
$ws.Range("G92").Value = 1308.1
$ws.Range("G93").Value = 26750.66
$ws.Range("G96").Value = 30102.66
$ws.Range("E99").Value = 11.9            # LABOR HOURS:

# --- FRC: 1C block (rows 126-134) ---
$ws.Range("G126").Value = 36.5
$ws.Range("E126").Value = "0.33 HRS"     # SPLICING LABOR hours label
$ws.Range("G128").Value = 36.5
$ws.Range("G129").Value = 342.7
$ws.Range("G131").Value = 342.7
$ws.Range("E134").Value = 0.33           # LABOR HOURS:

# --- FRC: 85C block (rows 161-169) ---
$ws.Range("G161").Value = 69.45
$ws.Range("E161").Value = "0.63 HRS"     # SPLICING LABOR hours label
$ws.Range("G163").Value = 69.45
$ws.Range("G164").Value = 875.35
$ws.Range("G166").Value = 875.35
$ws.Range("E169").Value = 0.63           # LABOR HOURS:

# --- FRC: 4C block (rows 195-204) ---
$ws.Range("G195").Value = 1315.9
$ws.Range("E195").Value = "11.97 HRS"    # SPLICING LABOR hours label
$ws.Range("G197").Value = 1315.9
$ws.Range("G198").Value = 15023.18
$ws.Range("G201").Value = 16588.18
$ws.Range("E204").Value = 11.97          # LABOR HOURS:

# --- FRC: 8645C block (rows 228-236) ---
$ws.Range("G228").Value = 13.67
$ws.Range("E228").Value = "0.12 HRS"     # SPLICING LABOR hours label
$ws.Range("G230").Value = 13.67
$ws.Range("G231").Value = 130.83
$ws.Range("G233").Value = 130.83
$ws.Range("E236").Value = 0.12           # LABOR HOURS:

# --- Part F billing summary (rows 246-258) ---
$ws.Range("B246").Value = 487646.66
$ws.Range("B248").Value = 487646.66
$ws.Range("C251").Value = 3138.776
$ws.Range("C253").Value = 254395.18
$ws.Range("C254").Value = 32789.3
$ws.Range("C255").Value = 487646.66
$ws.Range("C256").Value = 487646.66
$ws.Range("C258").Value = 298.34

# --- Part C detailed pricing summary (rows 282-284) ---
$ws.Range("C282").Value = 5541
$ws.Range("D282").Value = 264.048
$ws.Range("C284").Value = 37731
$ws.Range("D284").Value = 3138.32

# --- (1) ENGINEERING COST estimate (rows 303-304) ---
$ws.Range("C303").Value = 8.58
$ws.Range("B304").Value = 8.58

# --- COST ESTIMATE AS PER REQUEST total (row 356) ---
$ws.Range("B356").Value = 91179.58
